# Auto update Excel log
# Appends the latest sensor poll rows captured since the last export to each
# sheet's log table (ALERTS, PIR, Humidity, mmWave). Columns are always
# Date | Timestamp | Hour | Location | Value | Status, stored as plain text
# (matching the existing inlineStr cells), so date/time/percent-looking
# values are pinned to text format before the write to stop Excel's COM
# layer from auto-coercing them into real dates/numbers.

$wb = $excel.ActiveWorkbook

function Append-LogRows {
    param($ws, $startRow, $rows)
    $r = $startRow
    foreach ($row in $rows) {
        $ws.Range("A$r").NumberFormat = "@"
        $ws.Range("A$r").Value = $row[0]
        $ws.Range("B$r").Value = $row[1]
        $ws.Range("C$r").Value = $row[2]
        $ws.Range("D$r").Value = $row[3]
        if ($row[4] -like "*%") {
            $ws.Range("E$r").NumberFormat = "@"
        }
        $ws.Range("E$r").Value = $row[4]
        $ws.Range("F$r").Value = $row[5]
        $r = $r + 1
    }
}

$ws = $wb.Worksheets.Item("ALERTS")
$ALERTSRows = @(
    @("2026-01-30", "15:47:17", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "15:48:25", "15:00", "Living Room", "CRITICAL", "FALL_DETECTED"),
    @("2026-01-30", "15:48:29", "15:00", "Living Room", "CRITICAL", "FALL_DETECTED"),
    @("2026-01-30", "15:48:55", "15:00", "Living Room", "CRITICAL", "FALL_DETECTED")
)
Append-LogRows $ws 7 $ALERTSRows

$ws = $wb.Worksheets.Item("PIR")
$PIRRows = @(
    @("2026-01-30", "15:47:04", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:47:09", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:47:17", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:47:19", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:47:24", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:47:29", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:47:34", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:47:39", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:47:40", "15:00", "Living Room", "RECOVERY_DETECTION", "Inactive"),
    @("2026-01-30", "15:48:29", "15:00", "Living Room", "RECOVERY_DETECTION", "Inactive"),
    @("2026-01-30", "15:48:30", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:48:30", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:48:35", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:48:40", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:48:45", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:48:50", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:48:56", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:49:00", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:49:05", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:49:10", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:49:15", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:49:20", "15:00", "Bathroom", "No Motion", "Inactive")
)
Append-LogRows $ws 118 $PIRRows

$ws = $wb.Worksheets.Item("Humidity")
$HumidityRows = @(
    @("2026-01-30", "15:47:05", "15:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-30", "15:47:10", "15:00", "Bathroom", "86.9%", "Active"),
    @("2026-01-30", "15:47:17", "15:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-30", "15:47:20", "15:00", "Bathroom", "86.9%", "Active"),
    @("2026-01-30", "15:47:25", "15:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-30", "15:47:35", "15:00", "Bathroom", "86.3%", "Active"),
    @("2026-01-30", "15:48:30", "15:00", "Bathroom", "87.7%", "Active"),
    @("2026-01-30", "15:48:30", "15:00", "Bathroom", "87.7%", "Active"),
    @("2026-01-30", "15:48:35", "15:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "15:48:45", "15:00", "Bathroom", "87.7%", "Active"),
    @("2026-01-30", "15:48:50", "15:00", "Bathroom", "87.7%", "Active"),
    @("2026-01-30", "15:49:15", "15:00", "Bathroom", "86.8%", "Active")
)
Append-LogRows $ws 69 $HumidityRows

$ws = $wb.Worksheets.Item("mmWave")
$mmWaveRows = @(
    @("2026-01-30", "15:47:41", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:48:29", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:48:33", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
)
Append-LogRows $ws 14 $mmWaveRows

